$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- 1) "Quer que o evento deixe de ter data_inicio e data_fim na tabela eventos?"
#        -> remove spell/gram proofErr wrapping and normalize "data_inicio"/"data_fim"
#           to "data início" / "data fim".
$targetText1 = "Quer que o evento deixe de ter data_inicio e data_fim na tabela eventos?"
$p1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq $targetText1) { $p1 = $cand; break }
}
if ($p1 -ne $null) {
    $xml1 = '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t xml:space="preserve">Quer que o evento deixe de ter </w:t></w:r>' +
        '<w:r><w:t>data início</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> e </w:t></w:r>' +
        '<w:r><w:t>data fim</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> na tabela eventos?</w:t></w:r>' +
        '</w:p>'
    $p1.Range.InsertXML($xml1)
}

# --- 2) "Se confirmar, começamos pelo ALTER TABLE correto."
#        -> remove spell proofErr wrapping "TABLE" and turn it into the "TABELE" typo.
$targetText2 = "Se confirmar, começamos pelo ALTER TABLE correto."
$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq $targetText2) { $p2 = $cand; break }
}
if ($p2 -ne $null) {
    $xml2 = '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t xml:space="preserve">Se confirmar, começamos pelo </w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">ALTER </w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>TABELE</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> correto</w:t></w:r>' +
        '<w:r><w:t>.</w:t></w:r>' +
        '</w:p>'
    $p2.Range.InsertXML($xml2)
}

# --- 3) The trailing empty paragraph (<w:p/>) becomes a paragraph holding a single
#        space, followed by a brand-new paragraph with the "Sim, quero que..." reply.
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($lastP.Range.Text.TrimEnd() -eq "") {
    $xml3 = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
    $lastP.Range.InsertXML($xml3)

    $spaceP = $d.Paragraphs.Item($d.Paragraphs.Count)
    $spaceP.Range.InsertParagraphAfter()

    $newP = $d.Paragraphs.Item($d.Paragraphs.Count)
    $xml4 = '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t xml:space="preserve">Sim, quero que </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">o evento deixe de ter </w:t></w:r>' +
        '<w:r><w:t>data início</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> e </w:t></w:r>' +
        '<w:r><w:t>data fim</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> na tabela eventos</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> p</w:t></w:r>' +
        '<w:r><w:t>orque ele pode ter várias datas no ano</w:t></w:r>' +
        '<w:r><w:t>.</w:t></w:r>' +
        '</w:p>'
    $newP.Range.InsertXML($xml4)
}
